# Fruta / hortaliza, semanal
# Insert a new weekly record at row 23, pushing existing rows 23-32 down to 24-33.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(23).Insert()

$ws.Cells.Item(23, 1).Value = 4
$ws.Cells.Item(23, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(23, 3).Value = "Los Lagos"
$ws.Cells.Item(23, 4).Value = 44841
$ws.Cells.Item(23, 5).Value = 10
$ws.Cells.Item(23, 6).Value = 100112013
$ws.Cells.Item(23, 7).Value = "Alcachofa"
$ws.Cells.Item(23, 8).Value = "Madrigal"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 150
$ws.Cells.Item(23, 11).Value = 11000
$ws.Cells.Item(23, 12).Value = 12000
$ws.Cells.Item(23, 13).Value = 11533
$ws.Cells.Item(23, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(23, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(23, 16).Value = 384
$ws.Cells.Item(23, 17).Value = 30
$ws.Cells.Item(23, 18).Value = "Hortaliza"
